# Generate Report for Handoff
#
# Updates the localization-status report: rows move from "In Translation"
# to "Ready for handoff", and the associated timestamps are refreshed to
# reflect the new handoff generation time. The widened status/date columns
# follow from the longer "Ready for handoff" text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-09-06 21:19:39"
$overview.Range("E1").ColumnWidth = 16.33
$overview.Range("F1").ColumnWidth = 16.33

# --- zh-cn sheet --------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-09-06 21:19:34"
$zhcn.Range("C1").ColumnWidth = 16.33

# --- de-de sheet --------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-09-06 21:19:39"
$dede.Range("C1").ColumnWidth = 16.33
